$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "Price"
$ws.Range("B4").Value = "₹53,990"
